{"js": "// Task title correction: \"\u0413\u0430\u043b\u0430\u043a\u0442\u0438\u0447\u0435\u0441\u043a\u0430\u044f \u043f\u043e\u043b\u0438\u0446\u0438\u044f\" -> \"\u0413\u0430\u043b\u0430\u043a\u0442\u0438\u0447\u0435\u0441\u043a\u0430\u044f \u043d\u0430\u043b\u043e\u0433\u043e\u0432\u0430\u044f\"\n// (the author renamed the problem \"d\" from \"Galactic police\" to \"Galactic tax office\").\nconst body = context.document.body;\n\nconst results = body.search(\"\u043f\u043e\u043b\u0438\u0446\u0438\u044f\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replace only the second word (\"\u043f\u043e\u043b\u0438\u0446\u0438\u044f\" -> \"\u043d\u0430\u043b\u043e\u0433\u043e\u0432\u0430\u044f\") so the leading\n  // \"\u0413\u0430\u043b\u0430\u043a\u0442\u0438\u0447\u0435\u0441\u043a\u0430\u044f \" text/run is untouched, matching the author's edit\n  // (title \"\u0413\u0430\u043b\u0430\u043a\u0442\u0438\u0447\u0435\u0441\u043a\u0430\u044f \u043f\u043e\u043b\u0438\u0446\u0438\u044f\" -> \"\u0413\u0430\u043b\u0430\u043a\u0442\u0438\u0447\u0435\u0441\u043a\u0430\u044f \u043d\u0430\u043b\u043e\u0433\u043e\u0432\u0430\u044f\").\n  const found = results.items[0];\n  found.insertText(\"\u043d\u0430\u043b\u043e\u0433\u043e\u0432\u0430\u044f\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Task title correction: \"\u0413\u0430\u043b\u0430\u043a\u0442\u0438\u0447\u0435\u0441\u043a\u0430\u044f \u043f\u043e\u043b\u0438\u0446\u0438\u044f\" -> \"\u0413\u0430\u043b\u0430\u043a\u0442\u0438\u0447\u0435\u0441\u043a\u0430\u044f \u043d\u0430\u043b\u043e\u0433\u043e\u0432\u0430\u044f\"\n# (problem \"d\" was renamed from \"Galactic police\" to \"Galactic tax office\").\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$findText = \"\u043f\u043e\u043b\u0438\u0446\u0438\u044f\"\n$replaceText = \"\u043d\u0430\u043b\u043e\u0433\u043e\u0432\u0430\u044f\"\n$matchCase = $true\n$matchWholeWord = $true\n$matchWildcards = $false\n$matchSoundsLike = $false\n$matchAllWordForms = $false\n$forward = $true\n$wrap = 1          # wdFindContinue\n$format = $false\n$replace = 2       # wdReplaceAll\n\n$find.Execute($findText, $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, $replaceText, $replace) | Out-Null\n"}
